$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (Volume/Number and week-coverage dates)
$ws.Range("A8").Value = "Volume 30   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/13/2023  Through  2/19/2023"

$countFmt = $ws.Range("G19").NumberFormat
$pctFmt = $ws.Range("L22").NumberFormat
$countFmtBold = $ws.Range("F21").NumberFormat
$pctFmtBold = $ws.Range("E21").NumberFormat

$ws.Range("D15").NumberFormat = $countFmt
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = $pctFmt
$ws.Range("E15").Value = -100
$ws.Range("C23").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Paste($ws.Range("F15"))
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = -33.333333333333
$ws.Range("L15").NumberFormat = $pctFmt
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 6.25
$ws.Range("I16").Value = 25
$ws.Range("J16").Value = 25
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 25
$ws.Range("M16").Value = -48.979591836734
$ws.Range("N16").Value = -82.014388489208
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 46
$ws.Range("J17").Value = 39
$ws.Range("K17").Value = 17.948717948717
$ws.Range("L17").Value = 58.620689655172
$ws.Range("M17").Value = 170.588235294118
$ws.Range("N17").Value = 9.523809523809
$ws.Range("D18").NumberFormat = $countFmt
$ws.Range("D18").Value = 6
$ws.Range("E18").NumberFormat = $pctFmt
$ws.Range("E18").Value = -83.333333333333
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -10
$ws.Range("I18").Value = 18
$ws.Range("J18").Value = 17
$ws.Range("K18").Value = 5.882352941176
$ws.Range("L18").Value = 28.571428571428
$ws.Range("M18").Value = -40
$ws.Range("N18").Value = -90.374331550802
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -44.444444444444
$ws.Range("F19").Value = 47
$ws.Range("H19").Value = -33.802816901408
$ws.Range("I19").Value = 84
$ws.Range("J19").Value = 122
$ws.Range("K19").Value = -31.147540983606
$ws.Range("L19").Value = 55.555555555555
$ws.Range("M19").Value = 115.384615384615
$ws.Range("N19").Value = 6.32911392405
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 14
$ws.Range("E20").Value = -71.428571428571
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 41
$ws.Range("H20").Value = -41.463414634146
$ws.Range("I20").Value = 46
$ws.Range("J20").Value = 54
$ws.Range("K20").Value = -14.814814814814
$ws.Range("L20").Value = 130
$ws.Range("M20").Value = 15
$ws.Range("N20").Value = -90.927021696252
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 57
$ws.Range("E21").Value = -49.122807017543
$ws.Range("G21").Value = 165
$ws.Range("H21").Value = -26.060606060606
$ws.Range("I21").Value = 221
$ws.Range("J21").Value = 260
$ws.Range("K21").Value = -15
$ws.Range("L21").Value = 60.144927536231
$ws.Range("M21").Value = 23.463687150838
$ws.Range("N21").Value = -77.145811789038
$ws.Range("D22").NumberFormat = $countFmt
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = $pctFmt
$ws.Range("E22").Value = -100
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = 0
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 31.818181818181
$ws.Range("F24").Value = 101
$ws.Range("H24").Value = 8.602150537634
$ws.Range("I24").Value = 164
$ws.Range("J24").Value = 151
$ws.Range("K24").Value = 8.609271523178
$ws.Range("L24").Value = 70.833333333333
$ws.Range("M24").Value = 100
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 400
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 85
$ws.Range("I25").Value = 70
$ws.Range("J25").Value = 43
$ws.Range("K25").Value = 62.790697674418
$ws.Range("L25").Value = 32.075471698113
$ws.Range("M25").Value = 11.111111111111
$ws.Range("D26").NumberFormat = $countFmt
$ws.Range("D26").Value = 4
$ws.Range("E26").NumberFormat = $pctFmt
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -83.333333333333
$ws.Range("J26").Value = 7
$ws.Range("K26").Value = -57.142857142857
$ws.Range("L26").Value = 0
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 7
$ws.Range("H27").Value = 75
$ws.Range("I27").Value = 14
$ws.Range("J27").Value = 7
$ws.Range("L27").Value = 75
$ws.Range("C28").NumberFormat = $countFmt
$ws.Range("C28").Value = 2
$ws.Range("D28").NumberFormat = $countFmt
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = $pctFmt
$ws.Range("E28").Value = 100
$ws.Range("F28").NumberFormat = $countFmt
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 4
$ws.Range("J28").Value = 2
$ws.Range("L28").Value = 300
$ws.Range("M28").Value = 300
$ws.Range("N28").Value = -55.555555555555
$ws.Range("C29").NumberFormat = $countFmt
$ws.Range("C29").Value = 1
$ws.Range("D29").NumberFormat = $countFmt
$ws.Range("D29").Value = 1
$ws.Range("E29").NumberFormat = $pctFmt
$ws.Range("E29").Value = 0
$ws.Range("F29").NumberFormat = $countFmt
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 2
$ws.Range("J29").Value = 2
$ws.Range("L29").Value = 100
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = -77.777777777777

